$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# The blank spacer row (old row 13) is removed entirely; everything below
# (old row 14 onward) shifts up by one. This merges the taller/bordered
# spacer row (old row 14, ht=15.6 with styled C column) into the new row 13,
# and the whole data table (old rows 16-91) shifts up to rows 15-90.
$ws.Rows.Item(13).Delete()

# Re-select the (now) row 13 in full, matching the saved selection state.
$ws.Rows.Item(13).Select()
